$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force columns D and E (Price, Volume(1h)) to be treated as text,
# matching the original inlineStr cell type, so numeric-looking values
# (e.g. "0.150", "3.40", "0.0000328") are not coerced into numbers.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "71.832.13"
$ws.Range("E2").Value = "  +0.18%  "
$ws.Range("D3").Value = "4.003.12"
$ws.Range("E3").Value = "  -0.79%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").Value = "529.13"
$ws.Range("E5").Value = "  +0.38%  "
$ws.Range("D6").Value = "150.99"
$ws.Range("E6").Value = "  +1.74%  "
$ws.Range("D7").Value = "0.692"
$ws.Range("E7").Value = "  +10.00%  "
$ws.Range("E8").Value = "  -0.01%  "
$ws.Range("D9").Value = "0.743"
$ws.Range("E9").Value = "  +0.03%  "
$ws.Range("E10").Value = "  -3.80%  "
$ws.Range("D11").Value = "0.0000328"
$ws.Range("E11").Value = "  -4.70%  "
$ws.Range("D12").Value = "47.47"
$ws.Range("E12").Value = "  +1.38%  "
$ws.Range("D13").Value = "10.61"
$ws.Range("E13").Value = "  -2.90%  "
$ws.Range("D14").Value = "4.638.54"
$ws.Range("E14").Value = "  -0.76%  "
$ws.Range("D15").Value = "3.997.59"
$ws.Range("E15").Value = "  -0.75%  "
$ws.Range("D16").Value = "13.97"
$ws.Range("E16").Value = "  -2.57%  "
$ws.Range("D17").Value = "20.58"
$ws.Range("E17").Value = "  -4.04%  "
$ws.Range("E18").Value = "  -1.16%  "
$ws.Range("D19").Value = "1.18"
$ws.Range("E19").Value = "  -3.72%  "
$ws.Range("D20").Value = "71.699.97"
$ws.Range("E20").Value = "  -0.03%  "
$ws.Range("D21").Value = "426.57"
$ws.Range("E21").Value = "  -4.18%  "
$ws.Range("D22").Value = "97.33"
$ws.Range("E22").Value = "  +2.68%  "
$ws.Range("D23").Value = "3.47"
$ws.Range("E23").Value = "  -3.60%  "
$ws.Range("D24").Value = "4.17"
$ws.Range("E24").Value = "  +2.61%  "
$ws.Range("D25").Value = "14.36"
$ws.Range("E25").Value = "  -0.67%  "
$ws.Range("D26").Value = "11.23"
$ws.Range("E26").Value = "  -7.43%  "
$ws.Range("D27").Value = "10.69"
$ws.Range("E27").Value = "  -3.99%  "
$ws.Range("D28").Value = "5.83"
$ws.Range("E28").Value = "  +0.77%  "
$ws.Range("D29").Value = "36.54"
$ws.Range("E29").Value = "  -2.03%  "
$ws.Range("D30").Value = "3.58"
$ws.Range("E30").Value = "  +22.62%  "
$ws.Range("D31").Value = "13.35"
$ws.Range("E31").Value = "  -3.15%  "
$ws.Range("E32").Value = "  -2.30%  "
$ws.Range("D33").Value = "676.01"
$ws.Range("E33").Value = "  -3.35%  "
$ws.Range("D34").Value = "6.99"
$ws.Range("E34").Value = "  -0.34%  "
$ws.Range("B35").Value = "InjectiveProtocol"
$ws.Range("C35").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D35").Value = "43.98"
$ws.Range("E35").Value = "  +6.06%  "
$ws.Range("B36").Value = "OKB"
$ws.Range("C36").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D36").Value = "65.32"
$ws.Range("E36").Value = "  -3.95%  "
$ws.Range("D37").Value = "0.436"
$ws.Range("E37").Value = "  -2.56%  "
$ws.Range("D38").Value = "0.151"
$ws.Range("E38").Value = "  -1.65%  "
$ws.Range("D39").Value = "0.0₃0826"
$ws.Range("E39").Value = "  -9.79%  "
$ws.Range("E40").Value = "  -3.35%  "
$ws.Range("E41").Value = "  -0.02%  "
$ws.Range("D42").Value = "0.998"
$ws.Range("E42").Value = "  -0.21%  "
$ws.Range("D43").Value = "0.0485"
$ws.Range("E43").Value = "  -2.12%  "
$ws.Range("D44").Value = "3.16"
$ws.Range("E44").Value = "  +1.03%  "
$ws.Range("D45").Value = "0.150"
$ws.Range("E45").Value = "  +1.78%  "
$ws.Range("B46").Value = "Fetch.AI"
$ws.Range("C46").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D46").Value = "2.63"
$ws.Range("E46").Value = "  -9.01%  "
$ws.Range("B47").Value = "ApeXProtocol"
$ws.Range("C47").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D47").Value = "3.40"
$ws.Range("E47").Value = "  -4.15%  "
$ws.Range("D48").Value = "9.57"
$ws.Range("E48").Value = "  +2.47%  "
$ws.Range("D49").Value = "2.98"
$ws.Range("E49").Value = "  -5.96%  "
$ws.Range("D50").Value = "0.000271"
$ws.Range("E50").Value = "  -3.01%  "
$ws.Range("D51").Value = "146.10"
$ws.Range("E51").Value = "  +1.69%  "
